# Parts list for v2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the obsolete "Bluetooth Module" row (old row 9) ---
# This shifts the old row 10 (Piezo Speaker / Horn) up to row 9 and
# shrinks the table/dimension from F10 to F9.
$ws.Rows.Item(9).Delete()

# --- Row 4: Arduino Motorshield -> Arduino Motor Shield R3, with Purpose + URL ---
$ws.Range("F4").Value = "http://arduino.cc/en/Main/ArduinoMotorShieldR3"
$ws.Range("A4").Value = "Arduino Motor Shield R3 "

# --- Row 3: Arduino Uno -> Arduino Uno R3, with Purpose + URL ---
$ws.Range("A3").Value = "Arduino Uno R3"

$ws.Range("C4").Value = "Motor Shield"
$ws.Range("C3").Value = "Core"

$ws.Range("F3").Value = "http://www.arduino.cc/en/Main/arduinoBoardUno"

# --- Row 2: DF Robot chassis now has a product URL ---
$ws.Range("F2").Value = "http://www.dfrobot.com/index.php?route=product/product&product_id=97"

# --- Rows 5-7: LEDs now have a brand name and a price ---
$ws.Range("A5").Value = "Kingbright"
$ws.Range("A6").Value = "Kingbright Ultra Red"
$ws.Range("A7").Value = "Kingbright"

$ws.Range("D5").Value = 0.25
$ws.Range("D6").Value = 0.25
$ws.Range("D7").Value = 0.25

# --- Hyperlinks (relationships) for the newly added product URLs ---
$ws.Hyperlinks.Add($ws.Range("F4"), "http://arduino.cc/en/Main/ArduinoMotorShieldR3")
$ws.Hyperlinks.Add($ws.Range("F3"), "http://www.arduino.cc/en/Main/arduinoBoardUno")
$ws.Hyperlinks.Add($ws.Range("F2"), "http://www.dfrobot.com/index.php?route=product/product&product_id=97")

# Give the new price cells the same currency formatting used by D2.
$ws.Range("D2").Copy()
$ws.Range("D5:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the Hyperlink cell style (copied from the pre-existing F8
# hyperlink cell) to the newly created hyperlink cells, since adding a
# hyperlink through the object model otherwise generates its own style.
$ws.Range("F8").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Resize the table to match the new data extent (A1:F9) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F9"))

# --- Cosmetic: match the saved selection from the source file ---
$ws.Range("A5").Select()
